# Arreglos menores en visualizacion de saldo y nueva prueba de gestion de menu implementada
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) GestionarMenuControllerTest: add a new "Cupos negativos" test case row,
#    inserted right after the "Cadena vacia" row (old row 18), before the
#    "Cadena con solo numeros" row. Everything below shifts down by one row
#    and the "cupos" factor merge (A15:A18) grows to A15:A19.
# ---------------------------------------------------------------------------
$wsMenu = $wb.Worksheets.Item("GestionarMenuControllerTest")

# Insert a fresh row at 18 (pushes old rows 18-23 down to 19-24).
$wsMenu.Rows.Item(18).Insert()

# Clone the border/fill/font formatting of the row directly above (row 17,
# which is the last "middle" row of the cupos group) onto the new row so the
# A-column box border and table styling stay visually continuous.
$wsMenu.Range("A17:D17").Copy()
$wsMenu.Range("A18:D18").PasteSpecial(-4122)

# Fill in the new test case content.
$wsMenu.Range("B18").Value = "Cupos negativos"
$wsMenu.Range("C18").Value = '("10 febrero 2026", "Compota", "frutas", "-5")'
$wsMenu.Range("D18").Value = $false

# ---------------------------------------------------------------------------
# 2) RegisterControllerTest: the boolean-returning method now documents
#    returning true/false instead of 1/0.
# ---------------------------------------------------------------------------
$wsRegister = $wb.Worksheets.Item("RegisterControllerTest")

$newSpecText = "Los inputs deben ser distintos de null y vacío`nla password debe ser longitud >= 8, password==confirmPassword,  de ser asi retorna true`nCaso contrario retorna false"
$wsRegister.Range("C11").Value = $newSpecText
